# Redes-Sociales.xlsx update — adds a "Url" column (as the new column A) to every
# sheet, reformats the follower/like/etc. counts from abbreviated text
# ("18 mill.", "10,5 M", "124 k", ...) to plain grouped-digit text
# ("18.000.000", "105.000.000", "124.000", ...), and appends one new brand row
# to the Facebook, Instagram and Youtube sheets (Aristocrazy / Aristocrazy /
# MrBeast respectively).

function Looks-Numeric {
    param($val)
    if ($val -match '^-?[0-9]+$') { return $true }
    if ($val -match '^-?[0-9]+\.[0-9]+$') { return $true }
    return $false
}

# Writes $val into (row,col), forcing a Text-typed cell even when $val looks
# like a number (e.g. "115.000" or "66") so Excel doesn't silently coerce it
# into a numeric value.
function Set-Cell {
    param($ws, $row, $col, $val)
    if (Looks-Numeric $val) {
        $ws.Cells.Item($row, $col).NumberFormat = "@"
    }
    $ws.Cells.Item($row, $col).Value2 = $val
}

# Inserts a new column A ("Url") ahead of the existing header row, copying the
# header style (bold + border, style index 1 in this workbook) from the cell
# that used to be A1 (now B1) so the new header cell matches its neighbours.
function Add-UrlColumn {
    param($ws, $header)
    $ws.Columns.Item(1).Insert()
    $ws.Cells.Item(1, 2).Copy($ws.Cells.Item(1, 1))
    Set-Cell $ws 1 1 $header
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Facebook
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)
Add-UrlColumn $ws "Url"

Set-Cell $ws 2 1 "https://www.facebook.com/PandoraEspana/"
Set-Cell $ws 2 3 "18.000.000"
Set-Cell $ws 2 4 "18.000.000"

Set-Cell $ws 3 1 "https://www.facebook.com/tousjewelry"
Set-Cell $ws 3 3 "24.000.000"
Set-Cell $ws 3 4 "24.000.000"

Set-Cell $ws 4 1 "https://www.facebook.com/AristocrazySpain/"
Set-Cell $ws 4 2 "Aristocrazy"
Set-Cell $ws 4 3 "115.000"
Set-Cell $ws 4 4 "112.000"

# ---------------------------------------------------------------------------
# Sheet 2: Instagram
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)
Add-UrlColumn $ws "Url"

Set-Cell $ws 2 1 "https://www.instagram.com/theofficialpandora/"
Set-Cell $ws 2 3 "105.000.000"
Set-Cell $ws 2 5 "4094"

Set-Cell $ws 3 1 "https://www.instagram.com/tousjewelry/"
Set-Cell $ws 3 3 "19.000.000"

Set-Cell $ws 4 1 "https://www.instagram.com/aristocrazy/"
Set-Cell $ws 4 2 "aristocrazy"
Set-Cell $ws 4 3 "240.000"
Set-Cell $ws 4 4 "1"
Set-Cell $ws 4 5 "2975"

# ---------------------------------------------------------------------------
# Sheet 3: Youtube
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)
Add-UrlColumn $ws "Url"

Set-Cell $ws 2 1 "https://www.youtube.com/user/TheOfficialPandora"
Set-Cell $ws 2 4 "124.000"

Set-Cell $ws 3 1 "https://www.youtube.com/user/tousjewelry"
Set-Cell $ws 3 4 "26.600"

Set-Cell $ws 4 1 "https://www.youtube.com/@MrBeast"
Set-Cell $ws 4 2 "MrBeast"
Set-Cell $ws 4 3 "@MrBeast"
Set-Cell $ws 4 4 "242.000.000"
Set-Cell $ws 4 5 "780"

# ---------------------------------------------------------------------------
# Sheet 4: Twitter
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(4)
Add-UrlColumn $ws "Url"

Set-Cell $ws 2 1 "https://twitter.com/PANDORA_Corp"
Set-Cell $ws 2 3 "149.000"

Set-Cell $ws 3 1 "https://twitter.com/tousjewelry"
Set-Cell $ws 3 3 "795.000"
Set-Cell $ws 3 5 "134.000"

# ---------------------------------------------------------------------------
# Sheet 5: TikTok
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(5)
Add-UrlColumn $ws "Url"

Set-Cell $ws 2 1 "https://www.tiktok.com/@theofficialpandor"

Set-Cell $ws 3 1 "https://www.tiktok.com/@tousjewelry"
Set-Cell $ws 3 3 "268.800"
Set-Cell $ws 3 5 "1.400.000"

Write-Host "Redes-Sociales.xlsx updated: Url column added to all sheets, counts reformatted, new rows appended."
